$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regular rows: update Price (D) and Volume(1h) (E) columns only.
# D values that look numeric are prefixed with a leading quote so Excel
# keeps them as text (matching the source data which stores them as strings).
$priceVolumeUpdates = @(
    @{Row=2; D='48.188.43'; E='  +1.84%  '},
    @{Row=3; D='2.502.11'; E='  +0.46%  '},
    @{Row=4; D=$null; E='  -0.01%  '},
    @{Row=5; D='''321.44'; E='  -0.10%  '},
    @{Row=6; D='''108.07'; E='  -0.78%  '},
    @{Row=7; D='''0.526'; E='  +0.75%  '},
    @{Row=8; D='''1.00'; E='  -0.01%  '},
    @{Row=9; D='''0.540'; E='  +0.14%  '},
    @{Row=10; D='''39.88'; E='  +1.22%  '},
    @{Row=11; D='''20.20'; E='  +8.92%  '},
    @{Row=12; D='''0.0816'; E='  +0.79%  '},
    @{Row=13; D=$null; E='  +0.07%  '},
    @{Row=14; D='''7.16'; E='  -0.16%  '},
    @{Row=15; D='2.893.07'; E='  +0.47%  '},
    @{Row=16; D='2.501.87'; E='  +0.29%  '},
    @{Row=17; D='''0.844'; E='  -0.30%  '},
    @{Row=18; D='48.019.42'; E='  +1.64%  '},
    @{Row=19; D='''13.08'; E='  -2.42%  '},
    @{Row=20; D='''6.76'; E='  +1.88%  '},
    @{Row=21; D='0.0₃0946'; E='  +0.67%  '},
    @{Row=22; D='''2.78'; E='  +1.37%  '},
    @{Row=25; D=$null; E='  -0.37%  '},
    @{Row=27; D='''25.76'; E='  +0.21%  '},
    @{Row=28; D='''2.37'; E='  +3.11%  '},
    @{Row=29; D='''9.78'; E='  -1.92%  '},
    @{Row=30; D=$null; E='  +0.37%  '},
    @{Row=31; D='''35.26'; E='  +1.60%  '},
    @{Row=32; D='''49.24'; E='  -1.25%  '},
    @{Row=33; D='''19.55'; E='  -4.70%  '},
    @{Row=36; D='''0.0785'; E='  -0.05%  '},
    @{Row=37; D=$null; E='  -0.34%  '},
    @{Row=38; D='''4.64'; E='  -2.89%  '},
    @{Row=39; D=$null; E='  -0.33%  '},
    @{Row=40; D=$null; E='  -0.08%  '},
    @{Row=41; D='''121.11'; E='  +1.36%  '},
    @{Row=42; D=$null; E='  +0.61%  '},
    @{Row=43; D='''21.39'; E='  -7.14%  '},
    @{Row=44; D=$null; E='  +2.12%  '},
    @{Row=45; D='2.004.82'; E='  +0.45%  '},
    @{Row=46; D=$null; E='  +4.31%  '},
    @{Row=47; D=$null; E='  +3.98%  '},
    @{Row=48; D=$null; E='  -2.32%  '},
    @{Row=49; D=$null; E='  -1.71%  '},
    @{Row=50; D='''5.17'; E='  -0.21%  '},
    @{Row=51; D='''79.96'; E='  +3.12%  '}
)

foreach ($item in $priceVolumeUpdates) {
    $r = $item.Row
    if ($null -ne $item.D) {
        $ws.Cells.Item($r, 4).Value = $item.D
    }
    $ws.Cells.Item($r, 5).Value = $item.E
}

# Rows that were fully replaced (coin identity + price/volume swapped between adjacent ranks)
$fullRowUpdates = @(
    @{Row=23; B='BitcoinCash'; C='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D='''277.36'; E='  +12.34%  '},
    @{Row=24; B='Litecoin'; C='https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D='''71.95'; E='  +1.92%  '},
    @{Row=34; B='Filecoin'; C='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D='''5.34'; E='  +0.59%  '},
    @{Row=35; B='FirstDigitalUSD'; C='https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'; D='''1.00'; E='  -0.08%  '}
)

foreach ($item in $fullRowUpdates) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
}
